# Add 7 new rows (555-561) of landscaping data for 7/28/2025 (serial 45866),
# mirroring the existing shared "ABS(D-E)" formula in column F, and move the
# selection/active cell forward the same way the source workbook's author did.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 555
$lastNewRow  = 561

# Columns B..E, G..T for each new row (Date=A and Temp_Diff=F are handled separately below).
$rowsData = @(
    @("Flowering",    "Large",  72, 89, 0, 0.1, "Yes", 2, "Neutral", 9, 0.85, 73, 30.14, 6, 0.35, 4.7, 33, 20),
    @("Nonflowering", "Medium", 72, 89, 0, 0,   "Yes", 3, "Bright",  9, 0.85, 73, 30.14, 6, 0.35, 4.7, 33, 20),
    @("Nonflowering", "Small",  72, 89, 0, 0,   "Yes", 3, "Bright",  9, 0.85, 73, 30.14, 6, 0.35, 4.7, 33, 20),
    @("Nonflowering", "Medium", 72, 89, 0, 0.1, "Yes", 3, "Neutral", 9, 0.85, 73, 30.14, 6, 0.35, 4.7, 33, 20),
    @("Nonflowering", "Medium", 72, 89, 0, 0,   "Yes", 3, "Neutral", 9, 0.85, 73, 30.14, 6, 0.35, 4.7, 33, 20),
    @("Nonflowering", "Large",  72, 89, 0, 0.2, "Yes", 4, "Bright",  9, 0.85, 73, 30.14, 6, 0.35, 4.7, 33, 20),
    @("Tree",         "Medium", 72, 89, 0, 0.2, "Yes", 1, "Neutral", 9, 0.85, 73, 30.14, 6, 0.35, 4.7, 33, 20)
)

for ($i = 0; $i -lt $rowsData.Length; $i++) {
    $r = $firstNewRow + $i
    $vals = $rowsData[$i]

    # Column A: same date serial as the last day's rows before it (7/28/2025 = 45866).
    # Copy the date format (style) from the row directly above, then set the value,
    # so the new cells keep the existing "mm/dd/yyyy"-style numbering instead of a
    # brand new number format being invented.
    $ws.Range("A" + ($r - 1)).Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)
    $ws.Range("A" + $r).Value = 45866

    $arr = New-Object 'object[,]' 1,4
    $arr[0,0] = $vals[0]
    $arr[0,1] = $vals[1]
    $arr[0,2] = $vals[2]
    $arr[0,3] = $vals[3]
    $ws.Range("B" + $r + ":E" + $r).Value = $arr

    # Column F keeps the running shared "absolute temperature difference" formula.
    $ws.Range("F" + $r).Formula = "=ABS(D" + $r + "-E" + $r + ")"

    $arr2 = New-Object 'object[,]' 1,14
    for ($j = 0; $j -lt 14; $j++) {
        $arr2[0,$j] = $vals[4 + $j]
    }
    $ws.Range("G" + $r + ":T" + $r).Value = $arr2
}

$excel.Calculate()

# Move the view/selection the same amount the 7 new rows shifted everything (author
# had scrolled near the bottom of the sheet before adding the new day's data).
$ws.Range("I" + ($lastNewRow + 1)).Select()
